$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-12 from 2023-09-16 (45185)
# to 2023-10-05 (45204), keeping the existing date number format/style.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
